$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 3338462.2
$ws.Range("J17").Value = 3576820.2
$ws.Range("L17").Value = 10730460.6
$ws.Range("N17").Value = -10730796.6

$ws.Range("H62:N62").Value = @(3792.6428, 2922.7778, 5358.4, 2922.7778, 5358.4, -2298.7778, -6606.4)

$ws.Range("H65:N65").Value = @(3792.6428, 2922.7778, 5358.4, 14613.889, 26792, -11493.889, -33032)

$ws.Range("H112").Value = 1023.0465
$ws.Range("J112").Value = 1033.9269
$ws.Range("L112").Value = 3101.7807
$ws.Range("N112").Value = -5317.780699999999

$ws.Range("H129").Value = 118547.04
$ws.Range("J129").Value = 134320.39
$ws.Range("L129").Value = 402961.17
$ws.Range("N129").Value = -412961.17

$ws.Range("H132").Value = 3425.0688
$ws.Range("I132").Value = 3681.1538
$ws.Range("K132").Value = 11043.4614
$ws.Range("M132").Value = -8513.4614

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 2821.7646
$ws.Range("I61").Value = 2283.5715
$ws.Range("K61").Value = 2283.5715
$ws.Range("M61").Value = -2071.5715

$ws.Range("H74:N74").Value = @(19231902, 23809998, 3901.3, 23809998, 3901.3, -23809124, -5649.3)

$ws.Range("H77:N77").Value = @(19231902, 23809998, 3901.3, 119049990, 19506.5, -119045622, -28242.5)

$ws.Range("H102:N102").Value = @(1040.8334, 953.63635, 2000, 953.63635, 2000, 668.36365, -5244)

$ws.Range("H113").Value = 24759.2
$ws.Range("J113").Value = 24759.2
$ws.Range("L113").Value = 24759.2
$ws.Range("N113").Value = -33437.2

$ws.Range("H136").Value = 2821.7646
$ws.Range("I136").Value = 2283.5715
$ws.Range("K136").Value = 6850.7145
$ws.Range("M136").Value = -4300.7145

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H82:N82").Value = @(17370.666, 8378.5, 35355, 8378.5, 35355, -7995.5, -36121)

$ws.Range("H85:N85").Value = @(17370.666, 8378.5, 35355, 8378.5, 35355, -7052.5, -38007)

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 2000
$ws.Range("I16").Value = 0
$ws.Range("J16").Value = 2000
$ws.Range("K16").Value = 0
$ws.Range("L16").Value = 2000
$ws.Range("M16").ClearContents()
$ws.Range("N16").Value = -2574

$ws.Range("H22:N22").Value = @(182.2, 177.5, 201, 177.5, 201, 172.5, -901)

$ws.Range("H31:N31").Value = @(3115.6736, 1608.9667, 5494.684, 1608.9667, 5494.684, -1313.9667, -6084.684)

$ws.Range("H34:N34").Value = @(3115.6736, 1608.9667, 5494.684, 1608.9667, 5494.684, -1406.9667, -5898.684)

$ws.Range("H92").Value = 21734
$ws.Range("J92").Value = 21734
$ws.Range("L92").Value = 21734
$ws.Range("N92").Value = -26726

$ws.Range("H107:N107").Value = @(1961.625, 662.5714, 2972, 662.5714, 2972, 1257.4286, -6812)

$ws.Range("H113").Value = 2000
$ws.Range("I113").Value = 0
$ws.Range("J113").Value = 2000
$ws.Range("K113").Value = 0
$ws.Range("L113").Value = 2000
$ws.Range("M113").ClearContents()
$ws.Range("N113").Value = -6340

$ws.Range("H141").Value = 23340.035
$ws.Range("J141").Value = 23887.893
$ws.Range("L141").Value = 23887.893
$ws.Range("N141").Value = -34247.893

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H68").Value = 33884
$ws.Range("J68").Value = 50326.5
$ws.Range("L68").Value = 150979.5
$ws.Range("N68").Value = -152601.5

$ws.Range("H69:N69").Value = @(1899.4286, 1199.6666, 2424.25, 3598.9998, 7272.75, -2787.9998, -8894.75)

$ws.Range("H71").Value = 33884
$ws.Range("J71").Value = 50326.5
$ws.Range("L71").Value = 452938.5
$ws.Range("N71").Value = -461050.5

$ws.Range("H72:N72").Value = @(1899.4286, 1199.6666, 2424.25, 10796.9994, 21818.25, -6740.999400000001, -29930.25)

$ws.Range("H74").Value = 9926.888999999999
$ws.Range("J74").Value = 9926.888999999999
$ws.Range("L74").Value = 29780.667
$ws.Range("N74").Value = -31902.667

$ws.Range("H77").Value = 9926.888999999999
$ws.Range("J77").Value = 9926.888999999999
$ws.Range("L77").Value = 89342.00099999999
$ws.Range("N77").Value = -99950.00099999999

$ws.Range("H131").Value = 715.55
$ws.Range("J131").Value = 730.3723
$ws.Range("L131").Value = 2191.1169
$ws.Range("N131").Value = -12271.1169

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H33").Value = 7980
$ws.Range("J33").Value = 7980
$ws.Range("L33").Value = 7980
$ws.Range("N33").Value = -8484

$ws.Range("H49").Value = 4000
$ws.Range("J49").Value = 4000
$ws.Range("L49").Value = 4000
$ws.Range("N49").Value = -4368

$ws.Range("H70:N70").Value = @(2984747.8, 4327.727, 6263210, 4327.727, 6263210, -4057.727, -6263750)

$ws.Range("H73:N73").Value = @(2984747.8, 4327.727, 6263210, 4327.727, 6263210, -3391.727, -6265082)

$ws.Range("H126:N126").Value = @(4144.857, 4781.6665, 3667.25, 14344.9995, 11001.75, -11874.9995, -15941.75)

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H68").Value = 2500.4285
$ws.Range("I68").Value = 1751.5
$ws.Range("J68").Value = 2800
$ws.Range("K68").Value = 1751.5
$ws.Range("L68").Value = 2800
$ws.Range("M68").Value = -1002.5
$ws.Range("N68").Value = -4298

$ws.Range("H71").Value = 2500.4285
$ws.Range("I71").Value = 1751.5
$ws.Range("J71").Value = 2800
$ws.Range("K71").Value = 8757.5
$ws.Range("L71").Value = 14000
$ws.Range("M71").Value = -5013.5
$ws.Range("N71").Value = -21488

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H126:N126").Value = @(1806.8889, 1331.1818, 3900, 3993.5454, 11700, -1523.5454, -16640)

$ws.Range("H132").Value = 1610.6061
$ws.Range("I132").Value = 1066.08
$ws.Range("K132").Value = 3198.24
$ws.Range("M132").Value = -668.24
Write-Host "All edits applied."
